$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("B3").Value = "(0.053)"
$ws.Range("C3").Value = "(0.055)"

# Row 4
$ws.Range("B4").Value = "(0.055)"
$ws.Range("C4").Value = "(0.058)"

# Row 5
$ws.Range("B5").Value = "(0.061)"
$ws.Range("C5").Value = "(0.064)"

# Row 6
$ws.Range("B6").Value = "(0.067)"
$ws.Range("C6").Value = "(0.070)"

# Row 7
$ws.Range("B7").Value = "(0.072)"
$ws.Range("C7").Value = "(0.075)"

# Row 9
$ws.Range("B9").Value = "(0.043)"
$ws.Range("C9").Value = "(0.043)"

# Row 13
$ws.Range("B13").Value = "(0.032)"
$ws.Range("C13").Value = "(0.032)"

# Row 14
$ws.Range("B14").Value = "(0.086)"
$ws.Range("C14").Value = "(0.086)"

# Row 16
$ws.Range("C16").Value = "(0.057)"

# Row 17
$ws.Range("C17").Value = "(0.095)"

# Row 18
$ws.Range("C18").Value = "(0.028)"

# Row 19
$ws.Range("C19").Value = "(0.035)"

# Row 20
$ws.Range("C20").Value = "(0.024)"

# Row 21 (numeric observation counts)
$ws.Range("B21").Value = 8267
$ws.Range("C21").Value = 10204
